# Adiciona o popup de consulta com rolagem de tamanho (20 linhas extras de
# dados "Meta"/"Venda" = 5000/5000), converte a linha 8 para valores
# numéricos, e mantém a última linha (29) como texto "5000.0".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Linha 8: A8/B8 passam de texto ("1.2"/"0.0") para numérico (1.2 / 0)
$ws.Range("A8").Value = 1.2
$ws.Range("B8").Value = 0

# Linhas 9-28: novas linhas de rolagem do popup, numéricas, 5000/5000
for ($r = 9; $r -le 28; $r++) {
    $ws.Cells.Item($r, 1).Value = 5000
    $ws.Cells.Item($r, 2).Value = 5000
}

# Linha 29: última linha da rolagem, gravada como texto "5000.0" (colunas
# ainda não alinhadas, conforme mensagem do commit)
$ws.Range("A29").NumberFormat = "@"
$ws.Range("B29").NumberFormat = "@"
$ws.Range("A29").Value = "5000.0"
$ws.Range("B29").Value = "5000.0"
